$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 6.295968333333334
$ws.Range("H2").Value = 18.887905
$ws.Range("I2").Value = 0.5052862712055841
$ws.Range("J2").Value = 0.5052862712055841
$ws.Range("M2").Value = 110.642708
$ws.Range("N2").Value = 331.928124
$ws.Range("O2").Value = 0.5476418925386564
$ws.Range("P2").Value = 0.5476418925386564
$ws.Range("Q2").Value = 696.6029858822466
$ws.Range("R2").Value = 6269.426872940219
$ws.Range("S2").Value = 0.2767159298368269
$ws.Range("T2").Value = 0.2767159298368269
$ws.Range("G3").Value = 6.295968333333334
$ws.Range("H3").Value = 18.887905
$ws.Range("I3").Value = 0.5052862712055841
$ws.Range("J3").Value = 0.5052862712055841
$ws.Range("O3").Value = 0.3151072754333865
$ws.Range("P3").Value = 0.3151072754333865
$ws.Range("Q3").Value = 400.8178920034362
$ws.Range("R3").Value = 3607.361028030925
$ws.Range("S3").Value = 0.1592193802334868
$ws.Range("T3").Value = 0.1592193802334868
$ws.Range("G4").Value = 6.295968333333334
$ws.Range("H4").Value = 18.887905
$ws.Range("I4").Value = 0.5052862712055841
$ws.Range("J4").Value = 0.5052862712055841
$ws.Range("M4").Value = 27.72944133333333
$ws.Range("N4").Value = 83.18832399999999
$ws.Range("O4").Value = 0.1372508320279571
$ws.Range("P4").Value = 0.1372508320279571
$ws.Range("Q4").Value = 174.5836845356911
$ws.Range("R4").Value = 1571.25316082122
$ws.Range("S4").Value = 0.06935096113527041
$ws.Range("T4").Value = 0.06935096113527041
$ws.Range("I5").Value = 0.2025983155648483
$ws.Range("J5").Value = 0.2025983155648483
$ws.Range("M5").Value = 110.642708
$ws.Range("N5").Value = 331.928124
$ws.Range("O5").Value = 0.5476418925386564
$ws.Range("P5").Value = 0.5476418925386564
$ws.Range("Q5").Value = 279.3081854776253
$ws.Range("R5").Value = 2513.773669298628
$ws.Range("S5").Value = 0.1109513249610774
$ws.Range("T5").Value = 0.1109513249610775
$ws.Range("I6").Value = 0.2025983155648483
$ws.Range("J6").Value = 0.2025983155648483
$ws.Range("O6").Value = 0.3151072754333865
$ws.Range("P6").Value = 0.3151072754333865
$ws.Range("S6").Value = 0.06384020322503281
$ws.Range("T6").Value = 0.06384020322503282
$ws.Range("I7").Value = 0.2025983155648483
$ws.Range("J7").Value = 0.2025983155648483
$ws.Range("M7").Value = 27.72944133333333
$ws.Range("N7").Value = 83.18832399999999
$ws.Range("O7").Value = 0.1372508320279571
$ws.Range("P7").Value = 0.1372508320279571
$ws.Range("Q7").Value = 70.00063612978089
$ws.Range("R7").Value = 630.005725168028
$ws.Range("S7").Value = 0.02780678737873805
$ws.Range("T7").Value = 0.02780678737873805
$ws.Range("G8").Value = 3.639816666666666
$ws.Range("H8").Value = 10.91945
$ws.Range("I8").Value = 0.2921154132295675
$ws.Range("J8").Value = 0.2921154132295676
$ws.Range("M8").Value = 110.642708
$ws.Range("N8").Value = 331.928124
$ws.Range("O8").Value = 0.5476418925386564
$ws.Range("P8").Value = 0.5476418925386564
$ws.Range("Q8").Value = 402.7191726235333
$ws.Range("R8").Value = 3624.472553611799
$ws.Range("S8").Value = 0.159974637740752
$ws.Range("T8").Value = 0.159974637740752
$ws.Range("G9").Value = 3.639816666666666
$ws.Range("H9").Value = 10.91945
$ws.Range("I9").Value = 0.2921154132295675
$ws.Range("J9").Value = 0.2921154132295676
$ws.Range("O9").Value = 0.3151072754333865
$ws.Range("P9").Value = 0.3151072754333865
$ws.Range("Q9").Value = 231.7202956514722
$ws.Range("R9").Value = 2085.48266086325
$ws.Range("S9").Value = 0.09204769197486685
$ws.Range("T9").Value = 0.09204769197486687
$ws.Range("G10").Value = 3.639816666666666
$ws.Range("H10").Value = 10.91945
$ws.Range("I10").Value = 0.2921154132295675
$ws.Range("J10").Value = 0.2921154132295676
$ws.Range("M10").Value = 27.72944133333333
$ws.Range("N10").Value = 83.18832399999999
$ws.Range("O10").Value = 0.1372508320279571
$ws.Range("P10").Value = 0.1372508320279571
$ws.Range("Q10").Value = 100.9300827224222
$ws.Range("R10").Value = 908.3707445017999
$ws.Range("S10").Value = 0.04009308351394866
$ws.Range("T10").Value = 0.04009308351394866
